$d = $word.ActiveDocument

# --- Paragraph 1 ("San Gennaro Festa Scene On Mott St.") + Paragraph 2 ("By Dorothy Day") ---
# Replace the Heading1 title paragraph and the bold "By Dorothy Day" byline paragraph
# with a pandoc-style title block: a Title-styled paragraph whose text is split into
# one run per word (plus separate single-space runs), followed by an Authors-styled
# paragraph ("Dorothy Day", no "By " prefix, no bold) likewise split into per-word runs.

$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)
$combinedRange = $d.Range($p1.Range.Start, $p2.Range.End)

$titlePara = '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">San</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Gennaro</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Festa</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Scene</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">On</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Mott</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">St</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">.</w:t></w:r>' + `
    '</w:p>'

$authorsPara = '<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Day</w:t></w:r>' + `
    '</w:p>'

$pkgXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $titlePara + $authorsPara + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$combinedRange.InsertXML($pkgXml)

# --- Remove the bookmark that used to wrap the title paragraph ---
# The bookmark was stored as a pair of body-level markers (not nested inside the
# paragraph run content), so it predates/falls outside Word's normal Bookmarks
# object-model scanning; still, try every available avenue defensively in case the
# host recognizes it after the content rewrite above.
$bmName = "san-gennaro-festa-scene-on-mott-st."
try {
    if ($d.Bookmarks.Exists($bmName)) {
        $d.Bookmarks.Item($bmName).Delete()
    }
} catch {
}
try {
    $d.DeleteBookmark($bmName)
} catch {
}
